$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column cells keep their original text representation
# (values like "39.20" or "42.233.46" must not be auto-converted to numbers).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "42.233.46"
$ws.Range("E2").Value = "  +0.50%  "
$ws.Range("D3").Value = "2.287.40"
$ws.Range("E3").Value = "  -0.54%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "320.58"
$ws.Range("E5").Value = "  +1.20%  "
$ws.Range("D6").Value = "101.78"
$ws.Range("E6").Value = "  -2.11%  "
$ws.Range("D7").Value = "0.628"
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("E8").Value = "  +0.21%  "
$ws.Range("D9").Value = "0.605"
$ws.Range("E9").Value = "  -0.68%  "
$ws.Range("D10").Value = "39.20"
$ws.Range("E10").Value = "  -1.84%  "
$ws.Range("D11").Value = "0.0902"
$ws.Range("E11").Value = "  -0.96%  "
$ws.Range("D12").Value = "8.26"
$ws.Range("E12").Value = "  -2.26%  "
$ws.Range("E13").Value = "  -0.27%  "
$ws.Range("D14").Value = "0.959"
$ws.Range("E14").Value = "  -2.03%  "
$ws.Range("D15").Value = "15.11"
$ws.Range("E15").Value = "  -1.87%  "
$ws.Range("D16").Value = "2.632.50"
$ws.Range("E16").Value = "  -0.69%  "
$ws.Range("D17").Value = "2.274.04"
$ws.Range("E17").Value = "  -1.33%  "
$ws.Range("D18").Value = "42.328.04"
$ws.Range("E18").Value = "  +0.76%  "
$ws.Range("D19").Value = "7.36"
$ws.Range("E19").Value = "  -4.86%  "
$ws.Range("E20").Value = "  -0.26%  "
$ws.Range("D21").Value = "12.76"
$ws.Range("E21").Value = "  +28.96%  "
$ws.Range("D22").Value = "72.90"
$ws.Range("E22").Value = "  -0.33%  "
$ws.Range("D23").Value = "3.55"
$ws.Range("E23").Value = "  -0.79%  "
$ws.Range("D24").Value = "269.49"
$ws.Range("E24").Value = "  +4.04%  "
$ws.Range("E25").Value = "  -4.79%  "
$ws.Range("E26").Value = "  -0.28%  "
$ws.Range("D27").Value = "10.84"
$ws.Range("E27").Value = "  -1.34%  "
$ws.Range("D28").Value = "2.32"
$ws.Range("E28").Value = "  +2.27%  "
$ws.Range("D29").Value = "22.49"
$ws.Range("E29").Value = "  -1.77%  "
$ws.Range("D30").Value = "37.65"
$ws.Range("E30").Value = "  +4.27%  "
$ws.Range("D31").Value = "164.74"
$ws.Range("E31").Value = "  +0.15%  "
$ws.Range("D32").Value = "6.04"
$ws.Range("E32").Value = "  +2.48%  "
$ws.Range("D33").Value = "0.0873"
$ws.Range("E33").Value = "  -1.74%  "
$ws.Range("E34").Value = "  +0.91%  "
$ws.Range("E35").Value = "  -4.55%  "
$ws.Range("D36").Value = "2.49"
$ws.Range("E36").Value = "  -14.56%  "
$ws.Range("D37").Value = "4.57"
$ws.Range("E37").Value = "  -1.51%  "
$ws.Range("D38").Value = "0.0355"
$ws.Range("E38").Value = "  +0.68%  "
$ws.Range("D39").Value = "2.78"
$ws.Range("E39").Value = "  -5.07%  "
$ws.Range("D40").Value = "3.67"
$ws.Range("E40").Value = "  -0.29%  "
$ws.Range("D41").Value = "1.53"
$ws.Range("E41").Value = "  +3.48%  "
$ws.Range("E42").Value = "  +0.26%  "
$ws.Range("D43").Value = "68.16"
$ws.Range("E43").Value = "  -3.74%  "
$ws.Range("D44").Value = "0.224"
$ws.Range("E44").Value = "  -1.55%  "
$ws.Range("D45").Value = "91.29"
$ws.Range("E45").Value = "  -9.63%  "
$ws.Range("B46").Value = "Celestia"
$ws.Range("C46").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D46").Value = "12.17"
$ws.Range("E46").Value = "  -0.06%  "
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").Value = "114.65"
$ws.Range("E47").Value = "  +0.18%  "
$ws.Range("D48").Value = "79.29"
$ws.Range("E48").Value = "  +0.51%  "
$ws.Range("D49").Value = "8.94"
$ws.Range("E49").Value = "  -1.42%  "
$ws.Range("D50").Value = "1.608.18"
$ws.Range("E50").Value = "  +4.23%  "
$ws.Range("D51").Value = "5.22"
$ws.Range("E51").Value = "  -1.99%  "
